# fix several headers that weren't named in right format
$wb = $excel.ActiveWorkbook

# --- total_reg_and_cast ---
# Drop the unused "Registered Voters" column (always 0) and fix the
# precinct-column header casing.
$ws1 = $wb.Worksheets.Item("total_reg_and_cast")
$ws1.Activate()
$ws1.Columns.Item(2).Delete()
$ws1.Range("A1").Value = "precinct"
$ws1.Range("B2").Select()

# --- straightparty ---
$ws2 = $wb.Worksheets.Item("straightparty")
$ws2.Activate()
$ws2.Range("A57").Select()

# --- presidential ---
# Fix the precinct-column header casing.
$ws3 = $wb.Worksheets.Item("presidential")
$ws3.Activate()
$ws3.Range("A1").Value = "precinct"
$ws3.Range("A2").Select()

# --- ussenate ---
$ws4 = $wb.Worksheets.Item("ussenate")
$ws4.Activate()
$ws4.Range("C11").Select()

# --- cd03 ---
$ws5 = $wb.Worksheets.Item("cd03")
$ws5.Activate()
$ws5.Range("B4").Select()

# --- statehou87 ---
$ws6 = $wb.Worksheets.Item("statehou87")
$ws6.Activate()
$ws6.Range("D22").Select()

# Leave the workbook focused back on the first sheet, as before.
$ws1.Activate()
